# Daily attendance processing - reverse the order of names/emails in the
# "Recorded By" column (G) for each data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ", "
        $n = $parts.Count
        if ($n -gt 1) {
            $reversed = $parts[($n - 1)..0]
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
